$d = $word.ActiveDocument

# --- Step 1: modify the final existing paragraph ---
$lastPara = $d.Paragraphs.Last
$replaceRange = $lastPara.Range.Duplicate
$null = $replaceRange.Find.Execute('এলকিউজির অন্যতম প্রতিষ্ঠাতা কার্লো রোভেলি। ', $true, $false, $false, $false, $false, $true, 0, $false, 'এলকিউজির অন্যতম প্রধান প্রতিষ্ঠাতা কার্লো রোভেলি। তিনি কিছুদিন আগে বলেন', 2)

# Append the remaining 3 runs to complete the quoted sentence
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertAfter(', “')
$r.Font.Name = "SolaimanLipi"
$r.Font.NameBi = "SolaimanLipi"
$r.Collapse(0)
$r.InsertAfter('কোয়ান্টাম গ্র্যাভিটির অবস্থা বিশ বছর আগের চেয়ে অনেক ভাল। প্রতি দুই দিনে একদিন আমি এটা নিয়ে আশাবাদী থাকি।')
$r.Font.Name = "SolaimanLipi"
$r.Font.NameBi = "SolaimanLipi"
$r.Collapse(0)
$r.InsertAfter('"')
$r.Font.Name = "SolaimanLipi"
$r.Font.NameBi = "SolaimanLipi"

# --- New paragraph 0 ---
$prevLast = $d.Paragraphs.Last
$pr = $prevLast.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
# (spacer paragraph - leave empty)

# --- New paragraph 1 ---
$prevLast = $d.Paragraphs.Last
$pr = $prevLast.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$nr.Collapse(0)
$nr.InsertAfter('জনপ্রিয় বিজ্ঞানর পাঠকরা হয়তো লি স্মোলিনের কাছে এলকিউজি সম্পর্কে শুনেছেন। তিনি তত্ত্বটির আরেকজন প্রধান স্থপতি। ২০০০ সালে প্রকাশিত হয় তাঁর বই থ্রি রোডস টু কোয়ান্টাম গ্র্যাভিটি। পরবর্তীতে প্রকাশিত হয় আরেকটি বই দ্য ট্রাবল উইথ ফিজিক্স। এখানেও তিনি সংক্ষেপে এলকিউজি নিয়ে আলোচনা করেন। সম্প্রতি প্রকাশিত টাইম রিবর্ন বইয়েও তা করেন। রোভেলির সেভেন ব্রিফ লেসনস অন্য ফিজিক্স নামের বেস্ট-সেলিং বইয়েও এলকিউজির উল্লেখ আছে। সম্প্রতি প্রকাশিত রিয়েলিটি ইজ নট হোয়াট ইট সিমস বইয়েও তিনি এ আলোচনা করেছেন। ')
$wholeRange = $newPara.Range
$wholeRange.Font.Name = "SolaimanLipi"
$wholeRange.Font.NameBi = "SolaimanLipi"
$fr = $newPara.Range.Duplicate
$null = $fr.Find.Execute('থ্রি রোডস টু কোয়ান্টাম গ্র্যাভিটি', $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$fr.Italic = 1
$fr = $newPara.Range.Duplicate
$null = $fr.Find.Execute('দ্য ট্রাবল উইথ ফিজিক্স', $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$fr.Italic = 1
$fr = $newPara.Range.Duplicate
$null = $fr.Find.Execute('টাইম রিবর্ন', $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$fr.Italic = 1
$fr = $newPara.Range.Duplicate
$null = $fr.Find.Execute('সেভেন ব্রিফ লেসনস অন্য ফিজিক্স', $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$fr.Italic = 1
$fr = $newPara.Range.Duplicate
$null = $fr.Find.Execute('রিয়েলিটি ইজ নট হোয়াট ইট সিমস', $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$fr.Italic = 1

# --- New paragraph 2 ---
$prevLast = $d.Paragraphs.Last
$pr = $prevLast.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
# (spacer paragraph - leave empty)

# --- New paragraph 3 ---
$prevLast = $d.Paragraphs.Last
$pr = $prevLast.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$nr.Collapse(0)
$nr.InsertAfter('কোয়ান্টাম স্পেস বইটার উদ্দেশ্য হলো মানুষের ধারণায় ভারসাম্য তৈরি করা। আমি আপনাদের দেখাতে চাই, এলকিউজি শুধুই ভাল একটি তত্ত্ব নয়, এটি স্ট্রিং তত্ত্বের প্রকৃত ও বিশ্বাসযোগ্য একটি বিকল্প। কাজটা করতে গিয়ে আমি তত্ত্বটি সম্পর্কে স্মোলিন ও রোভেলি এ পর্যন্ত তাঁদের বইয়ে যা বেলছেন তার চেয়ে একটু বেশি তুলে ধরতে চাই। আমি আপনাদের বলতে চাই, এলকিউজি স্থান, কাল ও মহাবিশ্ব সম্পর্কে কী বলে। পাশাপাশি বলে দিতে চাই কেন ও কীভাবে তা বলে। ')
$wholeRange = $newPara.Range
$wholeRange.Font.Name = "SolaimanLipi"
$wholeRange.Font.NameBi = "SolaimanLipi"
$fr = $newPara.Range.Duplicate
$null = $fr.Find.Execute('কোয়ান্টাম স্পেস', $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$fr.Italic = 1

# --- New paragraph 4 ---
$prevLast = $d.Paragraphs.Last
$pr = $prevLast.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
# (spacer paragraph - leave empty)

# --- New paragraph 5 ---
$prevLast = $d.Paragraphs.Last
$pr = $prevLast.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$nr.Collapse(0)
$nr.InsertAfter('এ বইটি নিয়ে কাজ করতে ও লিখতে গিয়ে আমি স্মোলিন ও রোভেলি দুজনের কাছ থেকেই উল্লেখযোগ্য উৎসাহ, সমর্থন ও জ্ঞানের আলো পেয়ে ধন্য হয়েছি। এ বইটি আসলে তাঁদেরই গল্প। তবে আরও দুটি কথা খোলাখুলি বলে রাখি। বহু তাত্ত্বিকের বহু বছরের প্রচেষ্টার ফসল এলকিউজি। এ প্রচেষ্টাগুলো সম্পর্কে আমি সাধারণ মানুষের বোধগম্য করে যতটা সম্ভব বলে গিয়েছি। কারও অবদানের কথা সঠিকভাবে না উল্লেখ করা হলে বা উপেক্ষা করা হলে আমি আগেই ক্ষমা চেয়ে নিচ্ছি। এই বইটি মূলত তত্ত্বের প্রধান দুই ব্যক্তির কাজ নিয়ে লেখা। ফলে এলকিউজির নামে যত কাজ হয়েছে তার সবকিছুর সারমর্মও এতে পাওয়া যাবে না। ')
$wholeRange = $newPara.Range
$wholeRange.Font.Name = "SolaimanLipi"
$wholeRange.Font.NameBi = "SolaimanLipi"

# --- New paragraph 6 ---
$prevLast = $d.Paragraphs.Last
$pr = $prevLast.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
# (spacer paragraph - leave empty)

# --- New paragraph 7 ---
$prevLast = $d.Paragraphs.Last
$pr = $prevLast.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$nr = $newPara.Range
$nr.Collapse(0)
$nr.InsertAfter('বইটি তিন অংশে বিভক্ত। ')
$wholeRange = $newPara.Range
$wholeRange.Font.Name = "SolaimanLipi"
$wholeRange.Font.NameBi = "SolaimanLipi"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
